$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D column holds price text that looks numeric (e.g. "243.03"); a leading
# apostrophe forces Excel to store it as text instead of auto-converting it
# to a number, matching the source inline-string cell type.

$ws.Range("D2").Value = "30.489.73"
$ws.Range("D3").Value = "1.922.93"
$ws.Range("E3").Value = "  +1.66%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").Value = "'243.03"
$ws.Range("E5").Value = "  +1.19%  "
$ws.Range("E6").Value = "  +0.05%  "
$ws.Range("D7").Value = "'0.4704"
$ws.Range("E7").Value = "  -1.83%  "
$ws.Range("D8").Value = "'0.2877"
$ws.Range("E8").Value = "  -2.68%  "
$ws.Range("D9").Value = "'0.06753"
$ws.Range("E9").Value = "  +1.76%  "
$ws.Range("D10").Value = "'106.03"
$ws.Range("E10").Value = "  +5.06%  "
$ws.Range("D11").Value = "'18.27"
$ws.Range("E11").Value = "  -2.31%  "
$ws.Range("E12").Value = "  +1.77%  "
$ws.Range("D13").Value = "1.908.76"
$ws.Range("E13").Value = "  +1.00%  "
$ws.Range("D14").Value = "'5.302"
$ws.Range("E14").Value = "  +2.85%  "
$ws.Range("D15").Value = "'0.6588"
$ws.Range("E15").Value = "  -0.01%  "
$ws.Range("D16").Value = "'291.27"
$ws.Range("E16").Value = "  -5.77%  "
$ws.Range("D17").Value = "30.493.69"
$ws.Range("E17").Value = "  -0.94%  "
$ws.Range("E18").Value = "  +0.08%  "
$ws.Range("D19").Value = "'0.000007585"
$ws.Range("E19").Value = "  -0.65%  "
$ws.Range("D20").Value = "'12.92"
$ws.Range("E20").Value = "  -2.14%  "
$ws.Range("D21").Value = "2.149.97"
$ws.Range("E21").Value = "  +1.05%  "
$ws.Range("D22").Value = "'1.000"
$ws.Range("E22").Value = "  +0.11%  "
$ws.Range("D23").Value = "'5.265"
$ws.Range("E23").Value = "  +1.72%  "
$ws.Range("D24").Value = "'6.205"
$ws.Range("E24").Value = "  -0.13%  "
$ws.Range("D25").Value = "'9.365"
$ws.Range("E25").Value = "  +0.22%  "
$ws.Range("D26").Value = "'168.95"
$ws.Range("E26").Value = "  +0.76%  "
$ws.Range("D27").Value = "'21.41"
$ws.Range("E27").Value = "  +4.17%  "
$ws.Range("D28").Value = "'2.107"
$ws.Range("E28").Value = "  +7.81%  "
$ws.Range("D29").Value = "'0.1069"
$ws.Range("E29").Value = "  -5.97%  "
$ws.Range("D30").Value = "'1.367"
$ws.Range("E30").Value = "  +1.55%  "
$ws.Range("D31").Value = "'4.175"
$ws.Range("E31").Value = "  -0.75%  "
$ws.Range("D32").Value = "'3.984"
$ws.Range("E32").Value = "  -1.05%  "
$ws.Range("D33").Value = "'0.05027"
$ws.Range("E33").Value = "  -1.82%  "
$ws.Range("D34").Value = "'0.7419"
$ws.Range("E34").Value = "  -0.47%  "
$ws.Range("D35").Value = "'1.154"
$ws.Range("E35").Value = "  -1.01%  "
$ws.Range("D36").Value = "'0.02104"
$ws.Range("E36").Value = "  +6.80%  "
$ws.Range("D37").Value = "'2.728"
$ws.Range("E37").Value = "  +0.49%  "
$ws.Range("D38").Value = "'2.684"
$ws.Range("D39").Value = "'2.072"
$ws.Range("E39").Value = "  +0.61%  "
$ws.Range("D40").Value = "'110.15"
$ws.Range("E40").Value = "  +0.96%  "
$ws.Range("D41").Value = "'0.8725"
$ws.Range("E41").Value = "  -1.13%  "
$ws.Range("D42").Value = "'5.850"
$ws.Range("E42").Value = "  +3.31%  "
$ws.Range("D43").Value = "'0.4257"
$ws.Range("E43").Value = "  +0.81%  "
$ws.Range("E44").Value = "  +0.05%  "
$ws.Range("D45").Value = "'67.39"
$ws.Range("E45").Value = "  -0.86%  "
$ws.Range("D46").Value = "'49.02"
$ws.Range("E46").Value = "  +14.58%  "
$ws.Range("D47").Value = "'7.191"
$ws.Range("E47").Value = "  -2.66%  "
$ws.Range("D48").Value = "'9.277"
$ws.Range("E48").Value = "  +2.09%  "
$ws.Range("D49").Value = "'35.10"
$ws.Range("E49").Value = "  +0.24%  "
$ws.Range("D50").Value = "'0.1216"
$ws.Range("E50").Value = "  -1.53%  "
$ws.Range("D51").Value = "'0.2470"
$ws.Range("E51").Value = "  +9.90%  "
